$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.28
$ws.Range("G2").Value = 2.88
$ws.Range("H2").Value = 2.4
$ws.Range("I2").Value = 3.55
$ws.Range("J2").Value = 2.68
$ws.Range("K2").Value = 4.8
$ws.Range("L2").Value = 1.22
$ws.Range("N2").Value = 2.22
$ws.Range("P2").Value = 2.04
$ws.Range("Q2").Value = 1.56
$ws.Range("S2").Value = 2.3
$ws.Range("T2").Value = 1.4
$ws.Range("U2").Value = 1.98
$ws.Range("V2").Value = 1.39
$ws.Range("W2").Value = 1.55

# Row 3
$ws.Range("G3").Value = 1.63
$ws.Range("K3").Value = 5
$ws.Range("N3").Value = 5.9
$ws.Range("O3").Value = 1.18
$ws.Range("P3").Value = 2.72
$ws.Range("Q3").Value = 1.5
$ws.Range("S3").Value = 2.24
$ws.Range("T3").Value = 1.51
$ws.Range("V3").Value = 1.2
$ws.Range("W3").Value = 2.58
$ws.Range("AF3").Value = 12.5
$ws.Range("AN3").Value = 6.8
$ws.Range("AO3").Value = 46

# Row 4
$ws.Range("K4").Value = 4.2
$ws.Range("P4").Value = 2.16
$ws.Range("Q4").Value = 1.7

# Row 5
$ws.Range("F5").Value = 3.4
$ws.Range("H5").Value = 2.06
$ws.Range("K5").Value = 4.5
$ws.Range("R5").Value = 1.6
$ws.Range("T5").Value = 1.45
$ws.Range("U5").Value = 2.4
$ws.Range("AC5").Value = 14.5
$ws.Range("AH5").Value = 1000
$ws.Range("AK5").Value = 46
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# Row 6
$ws.Range("F6").Value = 7.4
$ws.Range("G6").Value = 8.2
$ws.Range("J6").Value = 4.9
$ws.Range("N6").Value = 5.2
$ws.Range("P6").Value = 2.44
$ws.Range("Q6").Value = 1.6
$ws.Range("T6").Value = 1.81
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 2.94
$ws.Range("W6").Value = 1.14
$ws.Range("Y6").Value = 11.5
$ws.Range("AA6").Value = 14.5
$ws.Range("AG6").Value = 27
$ws.Range("AH6").Value = 22
$ws.Range("AI6").Value = 38
$ws.Range("AJ6").Value = 230
$ws.Range("AO6").Value = 7
